# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures on the single worksheet, as produced by the scheduled
# "Updated symbol list" GitHub Actions run. Values are stored as plain text
# (the sheet uses inline/text strings throughout, e.g. "326.02", "-0.82%"),
# so each cell is written with a leading apostrophe to force text and the
# style is immediately reset to "Normal" so no numeric/quote-prefix
# formatting is introduced on cells that never had one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value.
$updates = [ordered]@{
    "D2"  = "326.02";  "E2"  = "-0.82%"
    "D3"  = "43.70";   "E3"  = "-2.02%"
    "D4"  = "5.555";   "E4"  = "0.26%"
    "D5"  = "0.08002"; "E5"  = "-1.09%"
                        "E6"  = "0.42%"
    "D8"  = "0.9418";  "E8"  = "-0.52%"
    "D9"  = "0.1152";  "E9"  = "-3.12%"
    "D10" = "0.1839";  "E10" = "-2.91%"
    "D11" = "11.91";   "E11" = "39.76%"
    "D12" = "0.09646"; "E12" = "-0.65%"
    "D13" = "0.04749"; "E13" = "16.85%"
    "D14" = "0.1068";  "E14" = "0.09%"
    "D15" = "0.001279";"E15" = "0.13%"
    "D16" = "0.04071"; "E16" = "-6.49%"
    "D17" = "0.005880";"E17" = "-1.62%"
    "D18" = "3.375"
    "D19" = "4.312"
                        "E20" = "-0.29%"
    "D21" = "0.1401";  "E21" = "2.27%"
    "D22" = "0.2513";  "E22" = "-2.90%"
    "D23" = "0.001246";"E23" = "0.63%"
    "D24" = "0.004293";"E24" = "-7.20%"
    "D25" = "0.0001193";"E25" = "-3.09%"
    "D26" = "0.0003756";"E26" = "-5.92%"
    "D38" = "0.02542"; "E38" = "-4.72%"
    "D39" = "0.05449"; "E39" = "-2.58%"
    "D40" = "0.007521";"E40" = "-2.08%"
    "D41" = "0.1385";  "E41" = "-0.82%"
    "D42" = "0.007389";"E42" = "-34.81%"
    "D43" = "0.002021";"E43" = "-4.21%"
    "D44" = "0.008318";"E44" = "-11.93%"
    "D45" = "0.00007119";"E45" = "0.56%"
    "D46" = "0.00000000752";"E46" = "0.21%"
                        "E47" = "1.33%"
                        "E48" = "39.68%"
    "D49" = "0.00002107";"E49" = "0.21%"
    "D50" = "0.0002006";"E50" = "0.21%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces the value to be stored as text rather than
    # being auto-converted to a number/percentage.
    $cell.Value = "'" + $updates[$addr]
    # Drop the implicit "quote prefix" formatting Excel applies so the cell
    # keeps the same (default) style it had before the edit.
    $cell.Style = "Normal"
}
